# Insert a new data row at row 125 (pushing the existing rows 125-162 down
# to 126-163) and populate it with the new "Poroto verde" price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(125).Insert()

$ws.Range("A125").Value = 10
$ws.Range("B125").Value = "Vega Modelo de Temuco"
$ws.Range("C125").Value = "La Araucanía"
$ws.Range("D125").Value = 44855
$ws.Range("E125").Value = 9
$ws.Range("F125").Value = 100112031
$ws.Range("G125").Value = "Poroto verde"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 1500
$ws.Range("K125").Value = 2000
$ws.Range("L125").Value = 2000
$ws.Range("M125").Value = 2000
$ws.Range("N125").Value = "$/kilo"
$ws.Range("O125").Value = "Provincia de Limarí"
$ws.Range("P125").Value = 2000
$ws.Range("Q125").Value = 1
$ws.Range("R125").Value = "Hortaliza"
